# LoginForm.xlsx revision:
#  - Sheet1 row 2: email changes to a "1" suffixed account, username shortened
#    from "Kukuh1998" to "Kukuh1"
#  - Sheet1 row 3: a brand-new row is populated (email/password/username),
#    with the email turned into a mailto: hyperlink just like row 2
#  - Sheet1 selection cursor moves from A7 to A4
#  - Sheet2 content is untouched (only cosmetic style bookkeeping shifts
#    there, which Excel manages internally)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Row 2: update the existing login's email + username -------------------
$ws1.Range("A2").Value = "kukuhpradipto1@gmail.com"
$ws1.Range("C2").Value = "Kukuh1"

# --- Row 3: fill in the new login row --------------------------------------
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "kukuhpradipto2@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("A3"), "mailto:kukuhpradipto2@gmail.com")
$ws1.Range("B3").Value = "Password123"
$ws1.Range("C3").Value = "Kukuh2"

# --- Move the saved cursor position to A4, like the author's session -------
$null = $ws1.Range("A4").Select()
